$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "202.173.124.126"
$ws.Range("B10").Value = 28.3621531
$ws.Range("C10").Value = 77.2828514
$ws.Range("D10").Value = 20
$ws.Range("E10").Value = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
$ws.Range("F10").Value = "Linux armv81"
$ws.Range("G10").Value = "2025-06-25T16:37:47.183Z"
